# [Fonds de solidarite] Add 2020-11-04 data
# Updates "nombre_aides" (col C) and "montant_total" (col D) for the rows
# whose underlying figures changed with the new data refresh. Values are
# written as text (matching the workbook's existing inline-string cell
# type for these columns) rather than as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 2; C = "1697"; D = "4240936.79" },
    @{ Row = 4; C = "1265"; D = "6531279.49" },
    @{ Row = 6; C = "875"; D = "4062290.54" },
    @{ Row = 15; C = "286"; D = "996376.49" },
    @{ Row = 17; C = "647"; D = "5193950.22" },
    @{ Row = 18; C = "186"; D = "1235955.93" },
    @{ Row = 19; C = "18"; D = "44471.00" },
    @{ Row = 21; C = "228"; D = "700413.14" },
    @{ Row = 23; C = "437"; D = "2641362.61" },
    @{ Row = 24; C = "197"; D = "1050686.27" },
    @{ Row = 31; C = "421"; D = "1228591.11" },
    @{ Row = 33; C = "790"; D = "5185067.39" },
    @{ Row = 35; C = "518"; D = "2779259.92" },
    @{ Row = 43; C = "401"; D = "1410679.40" },
    @{ Row = 45; C = "238"; D = "1109405.19" },
    @{ Row = 48; C = "742"; D = "2206051.35" },
    @{ Row = 49; C = "6"; D = "30176.00" },
    @{ Row = 50; C = "976"; D = "6156002.45" },
    @{ Row = 51; C = "707"; D = "3804250.37" },
    @{ Row = 54; C = "9733"; D = "26890502.14" },
    @{ Row = 55; C = "4"; D = "8100.00" },
    @{ Row = 57; C = "49"; D = "368009.00" },
    @{ Row = 58; C = "6644"; D = "33086824.23" },
    @{ Row = 59; C = "21"; D = "240000.00" },
    @{ Row = 60; C = "6500"; D = "27445116.36" },
    @{ Row = 61; C = "68"; D = "197070.00" },
    @{ Row = 62; C = "131"; D = "656607.46" },
    @{ Row = 74; C = "281"; D = "962768.30" },
    @{ Row = 75; C = "484"; D = "2953810.53" },
    @{ Row = 76; C = "291"; D = "2022164.29" },
    @{ Row = 79; C = "442"; D = "1310998.80" },
    @{ Row = 81; C = "1179"; D = "7698688.19" },
    @{ Row = 82; C = "632"; D = "3617090.36" },
    @{ Row = 85; C = "885"; D = "2415292.79" },
    @{ Row = 87; C = "6"; D = "22500.00" },
    @{ Row = 88; C = "1256"; D = "6630697.18" },
    @{ Row = 90; C = "896"; D = "4337866.84" },
    @{ Row = 92; C = "41"; D = "154152.23" },
    @{ Row = 100; C = "1320"; D = "3365754.28" }
)

foreach ($chg in $changes) {
    $row = $chg.Row

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $chg.C
    $cCell.ClearFormats()

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $chg.D
    $dCell.ClearFormats()
}
